$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.498.96"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "3.800.36"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "667.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").Value = "3.798.85"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "4.442.90"
$ws.Range("D16").Value = "3.801.99"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "70.494.51"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +20.28%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "474.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.45%  "
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D30").Value = "3.953.31"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  +6.55%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("E35").Value = "  +7.77%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "3.758.81"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.966"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +9.91%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  +5.26%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  +1.27%  "